$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.890.51'
$ws.Range('E2').Value = '  -0.32%  '
$ws.Range('D3').Value = '1.631.20'
$ws.Range('E3').Value = '  -0.68%  '
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.41'
$ws.Range('E5').Value = '  -0.68%  '
$ws.Range('E7').Value = '  -0.17%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '23.35'
$ws.Range('E8').Value = '  -0.91%  '
$ws.Range('E9').Value = '  -0.80%  '
$ws.Range('E10').Value = '  -0.51%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0880'
$ws.Range('E11').Value = '  -0.40%  '
$ws.Range('D12').Value = '1.861.22'
$ws.Range('E12').Value = '  -0.76%  '
$ws.Range('D13').Value = '1.628.99'
$ws.Range('E13').Value = '  -0.79%  '
$ws.Range('E14').Value = '  -1.35%  '
$ws.Range('E15').Value = '  -1.85%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.33'
$ws.Range('E16').Value = '  -0.37%  '
$ws.Range('D17').Value = '27.893.24'
$ws.Range('E17').Value = '  -0.31%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '229.30'
$ws.Range('E18').Value = '  -1.66%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.68'
$ws.Range('E19').Value = '  +1.03%  '
$ws.Range('E20').Value = '  -0.41%  '
$ws.Range('E21').Value = '  -0.18%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.33'
$ws.Range('E22').Value = '  -1.02%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.09'
$ws.Range('E23').Value = '  -4.89%  '
$ws.Range('E24').Value = '  -0.76%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '154.01'
$ws.Range('E25').Value = '  +0.83%  '
$ws.Range('E26').Value = '  -0.21%  '
$ws.Range('E27').Value = '  -0.34%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.49'
$ws.Range('E28').Value = '  -1.24%  '
$ws.Range('E29').Value = '  -0.15%  '
$ws.Range('E30').Value = '  -0.74%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.43'
$ws.Range('E32').Value = '  +0.56%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.12'
$ws.Range('E33').Value = '  +0.50%  '
$ws.Range('D34').Value = '1.390.04'
$ws.Range('E34').Value = '  -1.38%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.58'
$ws.Range('E35').Value = '  +0.00%  '
$ws.Range('E36').Value = '  +9.70%  '
$ws.Range('E37').Value = '  -0.76%  '
$ws.Range('E38').Value = '  +1.02%  '
$ws.Range('E39').Value = '  -1.55%  '
$ws.Range('E40').Value = '  -3.29%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.01'
$ws.Range('E41').Value = '  -1.24%  '
$ws.Range('E42').Value = '  -0.19%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.83'
$ws.Range('E43').Value = '  -2.17%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '65.63'
$ws.Range('E44').Value = '  -2.43%  '
$ws.Range('E45').Value = '  -1.72%  '
$ws.Range('D46').Value = '1.771.35'
$ws.Range('E46').Value = '  -0.74%  '
$ws.Range('E47').Value = '  -2.87%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '88.34'
$ws.Range('E48').Value = '  +0.26%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.102'
$ws.Range('E49').Value = '  +1.13%  '
$ws.Range('E50').Value = '  -0.49%  '
$ws.Range('E51').Value = '  +0.09%  '
